# Swap the match-data (columns B through AD) between paired rows.
# Column A (the running index) stays in place for each row position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 30  # column AD

$rowPairs = @(
    @(19, 20),
    @(38, 39),
    @(234, 235),
    @(293, 294)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rng1 = $ws.Range($ws.Cells.Item($r1, $firstCol), $ws.Cells.Item($r1, $lastCol))
    $rng2 = $ws.Range($ws.Cells.Item($r2, $firstCol), $ws.Cells.Item($r2, $lastCol))

    $vals1 = $rng1.Value2
    $vals2 = $rng2.Value2

    $rng1.Value2 = $vals2
    $rng2.Value2 = $vals1
}
